# Applies the FSReportSoftwareEngProject.docx edit:
#  1. Remove the title-page paragraph "DR. B.R. Ambedkar National Institute
#     of Technology, Jalandhar" (36pt) from its original spot (right after
#     the inserted picture, before the "Software Engineering Project" line).
#  2. Re-add the same institute name, but split across two centered 28pt
#     paragraphs ("...Technology, " / "Jalandhar"), inserted further down
#     the title page - right after the "CSE 5th Sem (G1)" block (after its
#     trailing empty paragraph) and before "Table of Contents".
#  3. Add a 397-twip (19.85pt) left indent to the 11 short skill/
#     contribution lines inside each team member's profile block
#     (Acquainted with / Novice at / Previous works / Contribution).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: remove the original "DR. B.R. Ambedkar ... Jalandhar" paragraph
# ---------------------------------------------------------------------
# Locate the paragraph and delete its whole range (including its paragraph
# mark) so the following paragraph slides up into its place.
# NOTE: Paragraph.Range.Text includes the trailing paragraph-mark (\r), so
# trim it off before comparing against the plain target string.
$titlePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd([char]13) -eq "DR. B.R. Ambedkar National Institute of Technology, Jalandhar") {
        $titlePara = $d.Paragraphs($i)
        break
    }
}
if ($titlePara -ne $null) {
    $titlePara.Range.Delete()
}

# ---------------------------------------------------------------------
# Step 2: insert the two new centered paragraphs after the empty
# paragraph that follows "CSE 5th Sem (G1)" (i.e. right before
# "Table of Contents").
# ---------------------------------------------------------------------
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd([char]13) -eq "Table of Contents") {
        $anchorIndex = $i - 1
        break
    }
}

if ($anchorIndex -gt 0) {
    $anchorPara = $d.Paragraphs($anchorIndex)
    $anchorPara.Range.InsertParagraphAfter()

    $para1 = $d.Paragraphs($anchorIndex + 1)
    $para1.Range.Text = "DR. B.R. Ambedkar National Institute of Technology, "
    $para1.Alignment = 1

    $para1.Range.InsertParagraphAfter()
    $para2 = $d.Paragraphs($anchorIndex + 2)
    $para2.Range.Text = "Jalandhar"
    $para2.Alignment = 1
}

# ---------------------------------------------------------------------
# Step 3: add a 397-twip (19.85pt) left indent to the skill/contribution
# lines for each of the three team members.
# ---------------------------------------------------------------------
$indentTexts = @(
    "Acquainted with: C/C++, Java, Android Development, API",
    "Novice at: Dart, Flutter, Firebase ",
    "Previous works: Photo Blog[1] ",
    "Contribution: User Authentication, and incorporating Google Map API",
    "Acquainted with: C/C++, Java, Dart, Flutter, Android Development, Firebase",
    "Novice at: APIs",
    "Previous works: Teentigada App[2]",
    "Contribution: Backend, and Connecting the app with server ",
    "Novice at: Java, APIs ",
    "Contribution: Frontend of the app"
)

# Note: "Acquainted with: C/C++, Java, Dart, Flutter, Android Development,
# Firebase" appears twice (2nd and 3rd team members); since we scan every
# paragraph in the document (not just the first hit) both occurrences are
# caught naturally.
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text.TrimEnd([char]13)
    foreach ($target in $indentTexts) {
        if ($t -eq $target) {
            $p.LeftIndent = 19.85
            break
        }
    }
}
